$wb = $excel.ActiveWorkbook

# Hunk 0: @@ -1171,22 +1171,22 @@ (sheet ALC)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("K11").Value = 66737.336
$ws.Range("I11").Value = 66737.336
$ws.Range("H11").Value = 66737.336
$ws.Range("M11").Value = -66597.336

# Hunk 1: @@ -2065,25 +2065,25 @@ (sheet ALC)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("L29").Value = 2400
$ws.Range("K29").Value = 903
$ws.Range("H29").Value = 550.5
$ws.Range("M29").Value = -622
$ws.Range("J29").Value = 800
$ws.Range("I29").Value = 301
$ws.Range("N29").Value = -2962

# Hunk 2: @@ -2512,25 +2512,25 @@ (sheet ALC)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("J38").Value = 426.66666
$ws.Range("N38").Value = -2023.99998
$ws.Range("H38").Value = 1697949.8
$ws.Range("L38").Value = 1279.99998

# Hunk 3: @@ -2772,25 +2772,25 @@ (sheet ALC)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("K43").Value = 1497.8
$ws.Range("M43").Value = -1428.8
$ws.Range("H43").Value = 1493.4166
$ws.Range("I43").Value = 1497.8
$ws.Range("J43").Value = 1490.2858
$ws.Range("L43").Value = 1490.2858
$ws.Range("N43").Value = -1628.2858

# Hunk 4: @@ -3522,25 +3522,25 @@ (sheet ALC)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("K58").Value = 9091434
$ws.Range("I58").Value = 3030478
$ws.Range("L58").Value = 8993.3334
$ws.Range("H58").Value = 1084240.8
$ws.Range("J58").Value = 2997.7778
$ws.Range("M58").Value = -9091284
$ws.Range("N58").Value = -9293.3334

# Hunk 5: @@ -7443,22 +7443,22 @@ (sheet ALC)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H136").Value = 49996.668
$ws.Range("N136").Value = -60196.668
$ws.Range("L136").Value = 49996.668
$ws.Range("J136").Value = 49996.668

# Hunk 6: @@ -7492,25 +7492,25 @@ (sheet ALC)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("L137").Value = 4523.549999999999
$ws.Range("K137").Value = 2964.2307
$ws.Range("I137").Value = 988.0769
$ws.Range("M137").Value = -414.2307000000001
$ws.Range("N137").Value = -9623.549999999999
$ws.Range("H137").Value = 1380.3585
$ws.Range("J137").Value = 1507.85

# Hunk 7: @@ -7544,25 +7544,25 @@ (sheet ALC)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("I138").Value = 4547
$ws.Range("J138").Value = 3160.152
$ws.Range("K138").Value = 13641
$ws.Range("N138").Value = -19760.456
$ws.Range("H138").Value = 3320.173
$ws.Range("L138").Value = 9480.456
$ws.Range("M138").Value = -8501

# Hunk 8: @@ -11374,25 +11374,25 @@ (sheet ARM)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("N74").Value = -4563.1875
$ws.Range("M74").Value = -652.7692999999999
$ws.Range("H74").Value = 2237.6206
$ws.Range("I74").Value = 1526.7693
$ws.Range("K74").Value = 1526.7693
$ws.Range("L74").Value = 2815.1875
$ws.Range("J74").Value = 2815.1875

# Hunk 9: @@ -11524,25 +11524,25 @@ (sheet ARM)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("J77").Value = 2815.1875
$ws.Range("H77").Value = 2237.6206
$ws.Range("I77").Value = 1526.7693
$ws.Range("K77").Value = 7633.8465
$ws.Range("M77").Value = -3265.8465
$ws.Range("L77").Value = 14075.9375
$ws.Range("N77").Value = -22811.9375

# Hunk 10: @@ -14195,25 +14195,25 @@ (sheet ARM)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("J132").Value = 3499.25
$ws.Range("K132").Value = 11485.1112
$ws.Range("H132").Value = 3785.9033
$ws.Range("M132").Value = -8955.111199999999
$ws.Range("L132").Value = 10497.75
$ws.Range("I132").Value = 3828.3704
$ws.Range("N132").Value = -15557.75

# Hunk 11: @@ -21232,22 +21232,22 @@ (sheet BSM)
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("M134").Value = -2403.3102
$ws.Range("K134").Value = 4938.3102
$ws.Range("H134").Value = 1773.0312
$ws.Range("I134").Value = 1646.1034

# Hunk 12: @@ -24471,25 +24471,25 @@ (sheet CRP)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("K58").Value = 1358.1177
$ws.Range("I58").Value = 1358.1177
$ws.Range("L58").Value = 52882.8
$ws.Range("H58").Value = 13068.272
$ws.Range("J58").Value = 52882.8
$ws.Range("M58").Value = -1155.1177
$ws.Range("N58").Value = -53288.8

# Hunk 13: @@ -28103,25 +28103,25 @@ (sheet CRP)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("J132").Value = 4264.7144
$ws.Range("K132").Value = 3923.7693
$ws.Range("H132").Value = 2342.8
$ws.Range("M132").Value = -1393.7693
$ws.Range("L132").Value = 12794.1432
$ws.Range("I132").Value = 1307.9231
$ws.Range("N132").Value = -17854.1432

# Hunk 14: @@ -28204,25 +28204,25 @@ (sheet CRP)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("L134").Value = 10028.4
$ws.Range("J134").Value = 3342.8
$ws.Range("M134").Value = -669.9231
$ws.Range("K134").Value = 3204.9231
$ws.Range("H134").Value = 1435.1613
$ws.Range("I134").Value = 1068.3077
$ws.Range("N134").Value = -15098.4

# Hunk 15: @@ -28305,25 +28305,25 @@ (sheet CRP)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H136").Value = 13068.272
$ws.Range("K136").Value = 4074.3531
$ws.Range("N136").Value = -163748.4
$ws.Range("L136").Value = 158648.4
$ws.Range("M136").Value = -1524.3531
$ws.Range("I136").Value = 1358.1177
$ws.Range("J136").Value = 52882.8

# Hunk 16: @@ -32002,25 +32002,25 @@ (sheet CUL)
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("N68").Value = -81205.47200000001
$ws.Range("K68").Value = 3021.16662
$ws.Range("H68").Value = 18607.586
$ws.Range("M68").Value = -2210.16662
$ws.Range("I68").Value = 1007.05554
$ws.Range("L68").Value = 79583.47200000001
$ws.Range("J68").Value = 26527.824

# Hunk 17: @@ -32155,25 +32155,25 @@ (sheet CUL)
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H71").Value = 18607.586
$ws.Range("N71").Value = -246862.416
$ws.Range("J71").Value = 26527.824
$ws.Range("K71").Value = 9063.49986
$ws.Range("M71").Value = -5007.49986
$ws.Range("L71").Value = 238750.416
$ws.Range("I71").Value = 1007.05554

# Hunk 18: @@ -35179,22 +35179,22 @@ (sheet CUL)
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 1482427.5
$ws.Range("J131").Value = 1482427.5
$ws.Range("N131").Value = -4457362.5
$ws.Range("L131").Value = 4447282.5

# Hunk 19: @@ -35228,25 +35228,25 @@ (sheet CUL)
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("J132").Value = 3034.2307
$ws.Range("K132").Value = 32820.0003
$ws.Range("H132").Value = 3227.6316
$ws.Range("M132").Value = -30290.0003
$ws.Range("L132").Value = 27308.0763
$ws.Range("I132").Value = 3646.6667
$ws.Range("N132").Value = -32368.0763

# Hunk 20: @@ -35488,25 +35488,25 @@ (sheet CUL)
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("L137").Value = 83342280
$ws.Range("K137").Value = 10357.5
$ws.Range("I137").Value = 3452.5
$ws.Range("M137").Value = -5257.5
$ws.Range("N137").Value = -83352480
$ws.Range("H137").Value = 20836432
$ws.Range("J137").Value = 27780760

# Hunk 21: @@ -39180,25 +39180,25 @@ (sheet GSM)
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("J70").Value = 4236.375
$ws.Range("L70").Value = 4236.375
$ws.Range("H70").Value = 87336.75
$ws.Range("M70").Value = -128616.94
$ws.Range("N70").Value = -4776.375
$ws.Range("I70").Value = 128886.94
$ws.Range("K70").Value = 128886.94

# Hunk 22: @@ -39324,25 +39324,25 @@ (sheet GSM)
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("N73").Value = -6108.375
$ws.Range("I73").Value = 128886.94
$ws.Range("J73").Value = 4236.375
$ws.Range("K73").Value = 128886.94
$ws.Range("M73").Value = -127950.94
$ws.Range("L73").Value = 4236.375
$ws.Range("H73").Value = 87336.75

# Hunk 23: @@ -42182,25 +42182,25 @@ (sheet GSM)
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("J132").Value = 3849.6667
$ws.Range("K132").Value = 8168.000100000001
$ws.Range("H132").Value = 3098.3333
$ws.Range("M132").Value = -5638.000100000001
$ws.Range("L132").Value = 11549.0001
$ws.Range("I132").Value = 2722.6667
$ws.Range("N132").Value = -16609.0001

# Hunk 24: @@ -42965,22 +42965,25 @@ (sheet LTW)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H6").Value = 22210
$ws.Range("M6").Value = -288
$ws.Range("L6").Value = 26572
$ws.Range("I6").Value = 400
$ws.Range("J6").Value = 26572
$ws.Range("K6").Value = 400
$ws.Range("N6").Value = -26796

# Hunk 25: @@ -47816,22 +47819,19 @@ (sheet LTW)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H105").Value = 0
$ws.Range("L105").Value = 0
$ws.Range("J105").Value = 0
$ws.Range("N105").ClearContents()

# Hunk 26: @@ -49136,25 +49136,25 @@ (sheet LTW)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("J132").Value = 3899.8
$ws.Range("K132").Value = 22296.4995
$ws.Range("H132").Value = 5826.5454
$ws.Range("M132").Value = -19766.4995
$ws.Range("L132").Value = 11699.4
$ws.Range("I132").Value = 7432.1665
$ws.Range("N132").Value = -16759.4

# Hunk 27: @@ -52310,22 +52310,22 @@ (sheet WVR)
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("L54").Value = 6922.125
$ws.Range("H54").Value = 6922.125
$ws.Range("J54").Value = 6922.125
$ws.Range("N54").Value = -7962.125

# Hunk 28: @@ -52708,22 +52708,22 @@ (sheet WVR)
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 10871490
$ws.Range("M62").Value = -14493796
$ws.Range("I62").Value = 14494420
$ws.Range("K62").Value = 14494420

# Hunk 29: @@ -52858,22 +52858,22 @@ (sheet WVR)
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H65").Value = 10871490
$ws.Range("K65").Value = 72472100
$ws.Range("I65").Value = 14494420
$ws.Range("M65").Value = -72468980

# Hunk 30: @@ -53627,25 +53627,25 @@ (sheet WVR)
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("N81").Value = -1003021
$ws.Range("K81").Value = 501039
$ws.Range("M81").Value = -499978
$ws.Range("I81").Value = 250519.5
$ws.Range("H81").Value = 333829.5
$ws.Range("J81").Value = 500449.5
$ws.Range("L81").Value = 1000899

# Hunk 31: @@ -53777,25 +53777,25 @@ (sheet WVR)
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("M84").Value = -2499891
$ws.Range("N84").Value = -5015103
$ws.Range("I84").Value = 250519.5
$ws.Range("K84").Value = 2505195
$ws.Range("H84").Value = 333829.5
$ws.Range("L84").Value = 5004495
$ws.Range("J84").Value = 500449.5

# Hunk 32: @@ -56132,25 +56132,25 @@ (sheet WVR)
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("J132").Value = 2423.5454
$ws.Range("K132").Value = 5748.6
$ws.Range("H132").Value = 2181.9524
$ws.Range("M132").Value = -3218.6
$ws.Range("L132").Value = 7270.6362
$ws.Range("I132").Value = 1916.2
$ws.Range("N132").Value = -12330.6362

# Hunk 33: @@ -56331,25 +56331,25 @@ (sheet WVR)
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 1720.3572
$ws.Range("K136").Value = 2700.75
$ws.Range("N136").Value = -11245.2
$ws.Range("L136").Value = 6145.200000000001
$ws.Range("M136").Value = -150.75
$ws.Range("I136").Value = 900.25
$ws.Range("J136").Value = 2048.4

